# Aula 5 update: reposition the picture on slide 4 ("Imagem 13") so it
# sits near the top of the slide instead of mid-slide.
#
# Original: <a:off x="7196362" y="2461450"/>
# New:      <a:off x="7196362" y="108548"/>
#
# PowerPoint COM exposes shape positions in points, while OOXML stores
# them in EMUs (1 point = 12700 EMU), so convert accordingly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$picture = $s.Shapes.Item(12)

$emuPerPoint = 12700
$newTopEmu = 108548

$picture.Top = $newTopEmu / $emuPerPoint
